$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new blank columns to interleave "standard deviation" columns
# after each of the existing data columns (motorActual_RPM, deltaP_inH20,
# flow_GPM, torque_lbft). Doing this left-to-right at fixed letters C, E, G, I
# works because each insertion pushes the remaining original columns one
# slot to the right, landing the next insertion point correctly.
$ws.Columns("C").Insert()
$ws.Columns("E").Insert()
$ws.Columns("G").Insert()
$ws.Columns("I").Insert()

# Header row for the new columns
$ws.Range("C1").Value = "motorActual_SD_RPM"
$ws.Range("E1").Value = "deltaP_SD_inH20"
$ws.Range("G1").Value = "flow_SD_GPM"
$ws.Range("I1").Value = "torque_SD_lbft"

# Re-apply the centered format (style index 3 in the original workbook, the
# same style column A already used) to the full header row.
$ws.Range("A1").Copy()
$ws.Range("B1:I1").PasteSpecial(-4122)

# motorActual_SD_RPM (column C)
$ws.Range("C2").Value = 25.242884740060848
$ws.Range("C3").Value = 23.172369322104313
$ws.Range("C4").Value = 20.97613215061341
$ws.Range("C5").Value = 23.066974660756863
$ws.Range("C6").Value = 17.380653612565894

# deltaP_SD_inH20 (column E)
$ws.Range("E2").Value = 1.1741379816699347
$ws.Range("E3").Value = 1.0373620390201124
$ws.Range("E4").Value = 0.27826246602799409
$ws.Range("E5").Value = 0.62319338892512044
$ws.Range("E6").Value = 0.91262259450444261

# flow_SD_GPM (column G)
$ws.Range("G2").Value = 0.03240370349203834
$ws.Range("G3").Value = 0.023021728866443807
$ws.Range("G4").Value = 0.078866976612522025
$ws.Range("G5").Value = 0.040373258476373686
$ws.Range("G6").Value = 0.10473776778220936

# torque_SD_lbft (column I)
$ws.Range("I2").Value = 0.0054772255750515442
$ws.Range("I3").Value = 0.0044721359549994844
$ws.Range("I4").Value = 0.0070710678118654034
$ws.Range("I5").Value = 0.0044721359549995832
$ws.Range("I6").Value = 0.0054772255750516656

# Apply the centered format (matching column A's style) to the new SD data
# columns only -- the original data columns (B, D, F, H) keep their own
# pre-existing styles untouched.
$ws.Range("A2").Copy()
$ws.Range("C2:C6").PasteSpecial(-4122)
$ws.Range("E2:E6").PasteSpecial(-4122)
$ws.Range("G2:G6").PasteSpecial(-4122)
$ws.Range("I2:I6").PasteSpecial(-4122)

$ws.Range("F14").Select()
